$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-converted to a number
# by the Excel input parser; force them to remain plain text first.
$textCells = @("D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D19", "D23", "D24", "D26", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.291.04"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.870.06"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4701"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").Value = "0.2873"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "0.06587"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "21.85"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").Value = "0.08019"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "97.24"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").Value = "1.872.20"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "5.128"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "0.6859"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "269.28"
$ws.Range("E16").Value = "  -3.40%  "
$ws.Range("D17").Value = "30.274.74"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "14.06"
$ws.Range("E18").Value = "  +3.69%  "
$ws.Range("D19").Value = "0.000007675"
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("D21").Value = "2.119.39"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "5.277"
$ws.Range("E23").Value = "  -1.68%  "
$ws.Range("D24").Value = "6.222"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").Value = "168.22"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "1.951"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "0.09868"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "4.383"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Value = "1.464"
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "4.082"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("D34").Value = "0.04709"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "0.7012"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +1.02%  "
$ws.Range("D39").Value = "2.624"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "6.306"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "72.30"
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("D42").Value = "1.955"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "0.8438"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "103.07"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.067"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.174"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "925.42"
$ws.Range("E49").Value = "  -5.96%  "
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").Value = "0.05678"
$ws.Range("E51").Value = "  +0.63%  "
